# Commit: "Added favicon and updated resume"
$d = $word.ActiveDocument

# 1. Contact line: update portfolio URL "ryanshafi.com" -> "https://r-s-hafi.github.io/"
#    (scoped to paragraph 2 so the later "Portfolio Website (ryanshafi.com)" line is untouched)
#    A temporary bookmark is dropped immediately before the target run so that the
#    text replacement doesn't get coalesced into the preceding (identically
#    formatted) single-space run; the bookmark is removed again right after, which
#    leaves only the intended single-run text change behind.
$p1 = $d.Paragraphs.Item(2)
$r1 = $p1.Range
$r1.Find.Execute("ryanshafi.com")
$barrier = $d.Bookmarks.Add("tmpBarrier", $d.Range($r1.Start, $r1.Start))
$r1.Text = "https://r-s-hafi.github.io/"
$d.Bookmarks("tmpBarrier").Delete()

# 2. Co-op bullet: remove the proofErr-induced run split around "as-built" by
#    re-asserting the full sentence as one Find/Replace over that paragraph.
$p2 = $d.Paragraphs.Item(23)
$text2 = "Created a new piping and instrumentation diagram (P&ID) as-built for pump house supplying fire and cooling water to site using BricsCAD software"
$p2.Range.Find.Execute($text2, $true, $false, $false, $false, $false, `
                        $true, 1, $false, $text2, 2)

# 3. Skills bullet: remove proofErr-induced run splits around "BricsCAD" and "Seeq"
$p3 = $d.Paragraphs.Item(43)
$text3 = "Engineering Software: APEN Plus, AutoCAD, BricsCAD, SolidWorks, Siemens NX, MS 365, Seeq, GMARS"
$p3.Range.Find.Execute($text3, $true, $false, $false, $false, $false, `
                        $true, 1, $false, $text3, 2)

# 4. Skills bullet: remove proofErr-induced run split around "yfinance"
$p4 = $d.Paragraphs.Item(45)
$text4 = "Programming: Python, HTML/CSS, JavaScript (learning) | PyQt5, matplotlib | APIs (yfinance, OpenAI) | Git/GitHub"
$p4.Range.Find.Execute($text4, $true, $false, $false, $false, $false, `
                        $true, 1, $false, $text4, 2)
